# "correção nos dados e inicio da analise PNAD 2009"
#
# The original sheet had an extra header-only row 6
# ("grandes regiões e unidades da federação") with no data, which pushed
# every region/state row (norte, rondônia, ..., goiás) one row lower than
# it should be. The fix removes that stray row entirely: Excel shifts all
# rows below it up by one, so "norte" (and its data) lands on row 6,
# "rondônia" on row 7, ..., and "goiás" ends on row 36 (the sheet shrinks
# from A1:G37 to A1:G36). No cell values actually change - they just move
# up one row - and the now-unused shared string for the removed header is
# dropped from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(6).Delete()
